$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the closing entry for 1-DIC-2023 ("cierre 5 DIC 23")
$ws.Range("B61").Value = 45261
$ws.Range("C61").Value = "a comprobar"
$ws.Range("D61").Value = 500

# Move the active selection to D62, matching the author's final cursor position
$ws.Range("D62").Select()
